# Append: 2025-09-17 01:41 JST
# Update the "取得日時" (retrieved datetime) column (A) for all existing
# data rows (2-19) on the active sheet ("ランサーズ") from
# "2025-09-17 01:13:07" to "2025-09-17 01:41:13". The values are stored as
# plain text, so they are written as strings to avoid Excel auto-converting
# them into date serials.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-17 01:41:13"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -ne $null -and $cell.Value -ne "") {
        $cell.Value = $newTimestamp
    }
}
